$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4090.7144
$ws.Range("I17").Value = 1300
$ws.Range("J17").Value = 4158.7803
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 12476.3409
$ws.Range("M17").Value = -3732
$ws.Range("N17").Value = -12812.3409
$ws.Range("H98").Value = 1249.7084
$ws.Range("I98").Value = 863.3182
$ws.Range("K98").Value = 863.3182
$ws.Range("M98").Value = 634.6818
$ws.Range("H100").Value = 4565.5
$ws.Range("I100").Value = 4900.6
$ws.Range("J100").Value = 2890
$ws.Range("K100").Value = 4900.6
$ws.Range("L100").Value = 2890
$ws.Range("M100").Value = -4359.6
$ws.Range("N100").Value = -3972
$ws.Range("H113").Value = 6377.294
$ws.Range("I113").Value = 6573.1665
$ws.Range("J113").Value = 6270.4546
$ws.Range("K113").Value = 6573.1665
$ws.Range("L113").Value = 6270.4546
$ws.Range("M113").Value = -3319.1665
$ws.Range("N113").Value = -12778.4546
$ws.Range("H122").Value = 1249.7084
$ws.Range("I122").Value = 863.3182
$ws.Range("K122").Value = 2589.9546
$ws.Range("M122").Value = -139.9546
$ws.Range("H138").Value = 2454.16
$ws.Range("I138").Value = 1502.5
$ws.Range("J138").Value = 2788.527
$ws.Range("K138").Value = 4507.5
$ws.Range("L138").Value = 8365.581
$ws.Range("M138").Value = 632.5
$ws.Range("N138").Value = -18645.581

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4253.299
$ws.Range("I32").Value = 2745.138
$ws.Range("J32").Value = 8857.157999999999
$ws.Range("K32").Value = 2745.138
$ws.Range("L32").Value = 8857.157999999999
$ws.Range("M32").Value = -2458.138
$ws.Range("N32").Value = -9431.157999999999
$ws.Range("H97").Value = 10481.689
$ws.Range("I97").Value = 9202.714
$ws.Range("J97").Value = 13839
$ws.Range("K97").Value = 9202.714
$ws.Range("L97").Value = 13839
$ws.Range("M97").Value = -8706.714
$ws.Range("N97").Value = -14831
$ws.Range("H122").Value = 3203
$ws.Range("I122").Value = 2625.7368
$ws.Range("K122").Value = 7877.2104
$ws.Range("M122").Value = -5427.2104

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 14051
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H64").Value = 2239
$ws.Range("I64").Value = 2065
$ws.Range("K64").Value = 2065
$ws.Range("M64").Value = -1840
$ws.Range("H67").Value = 2239
$ws.Range("I67").Value = 2065
$ws.Range("K67").Value = 2065
$ws.Range("M67").Value = -1285
$ws.Range("H99").Value = 75120.57000000001
$ws.Range("I99").Value = 202519
$ws.Range("J99").Value = 4343.6665
$ws.Range("K99").Value = 202519
$ws.Range("L99").Value = 4343.6665
$ws.Range("M99").Value = -201021
$ws.Range("N99").Value = -7339.6665

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3066.8438
$ws.Range("I31").Value = 2175.9546
$ws.Range("J31").Value = 5026.8
$ws.Range("K31").Value = 2175.9546
$ws.Range("L31").Value = 5026.8
$ws.Range("M31").Value = -1880.9546
$ws.Range("N31").Value = -5616.8
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 3066.8438
$ws.Range("I34").Value = 2175.9546
$ws.Range("J34").Value = 5026.8
$ws.Range("K34").Value = 2175.9546
$ws.Range("L34").Value = 5026.8
$ws.Range("M34").Value = -1973.9546
$ws.Range("N34").Value = -5430.8
$ws.Range("H58").Value = 4268.7915
$ws.Range("I58").Value = 4253.5264
$ws.Range("J58").Value = 4326.8
$ws.Range("K58").Value = 4253.5264
$ws.Range("L58").Value = 4326.8
$ws.Range("M58").Value = -4050.5264
$ws.Range("N58").Value = -4732.8
$ws.Range("H136").Value = 4268.7915
$ws.Range("I136").Value = 4253.5264
$ws.Range("J136").Value = 4326.8
$ws.Range("K136").Value = 12760.5792
$ws.Range("L136").Value = 12980.4
$ws.Range("M136").Value = -10210.5792
$ws.Range("N136").Value = -18080.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 405.5
$ws.Range("I2").Value = 196.28572
$ws.Range("J2").Value = 475.2381
$ws.Range("K2").Value = 1177.71432
$ws.Range("L2").Value = 2851.4286
$ws.Range("M2").Value = -1064.71432
$ws.Range("N2").Value = -3077.4286
$ws.Range("H9").Value = 333933.34
$ws.Range("I9").Value = 1000000
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 3000000
$ws.Range("L9").Value = 2700
$ws.Range("M9").Value = -2999776
$ws.Range("N9").Value = -3148
$ws.Range("H46").Value = 334766.2
$ws.Range("I46").Value = 1667058.5
$ws.Range("K46").Value = 5001175.5
$ws.Range("M46").Value = -5001084.5
$ws.Range("H107").Value = 332.75
$ws.Range("I107").Value = 194.8
$ws.Range("J107").Value = 677.625
$ws.Range("K107").Value = 584.4000000000001
$ws.Range("L107").Value = 2032.875
$ws.Range("M107").Value = 1335.6
$ws.Range("N107").Value = -5872.875
$ws.Range("H132").Value = 2927.4443
$ws.Range("I132").Value = 2800.8333
$ws.Range("J132").Value = 2990.75
$ws.Range("K132").Value = 25207.4997
$ws.Range("L132").Value = 26916.75
$ws.Range("M132").Value = -22677.4997
$ws.Range("N132").Value = -31976.75
$ws.Range("H137").Value = 2702.2727
$ws.Range("I137").Value = 1661.1428
$ws.Range("J137").Value = 4524.25
$ws.Range("K137").Value = 4983.428400000001
$ws.Range("L137").Value = 13572.75
$ws.Range("M137").Value = 116.5715999999993
$ws.Range("N137").Value = -23772.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 19060
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19060
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19060
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -19770
$ws.Range("H122").Value = 107986.445
$ws.Range("I122").Value = 157810.92
$ws.Range("K122").Value = 473432.76
$ws.Range("M122").Value = -470982.76

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8531.125
$ws.Range("I46").Value = 7800
$ws.Range("J46").Value = 9749.666999999999
$ws.Range("K46").Value = 7800
$ws.Range("L46").Value = 9749.666999999999
$ws.Range("M46").Value = -7612
$ws.Range("N46").Value = -10125.667
$ws.Range("H61").Value = 2549.7273
$ws.Range("I61").Value = 2725.7
$ws.Range("J61").Value = 790
$ws.Range("K61").Value = 2725.7
$ws.Range("L61").Value = 790
$ws.Range("M61").Value = -2523.7
$ws.Range("N61").Value = -1194
$ws.Range("H93").Value = 5176.1665
$ws.Range("I93").Value = 5651.6
$ws.Range("K93").Value = 5651.6
$ws.Range("M93").Value = -4403.6
$ws.Range("H113").Value = 2549.7273
$ws.Range("I113").Value = 2725.7
$ws.Range("J113").Value = 790
$ws.Range("K113").Value = 2725.7
$ws.Range("L113").Value = 790
$ws.Range("M113").Value = -555.6999999999998
$ws.Range("N113").Value = -5130
$ws.Range("H140").Value = 100429
$ws.Range("J140").Value = 100429
$ws.Range("L140").Value = 100429
$ws.Range("N140").Value = -110789

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 1630.7858
$ws.Range("I100").Value = 1115.6666
$ws.Range("K100").Value = 2231.3332
$ws.Range("M100").Value = -1690.3332
$ws.Range("H126").Value = 1730.5186
$ws.Range("I126").Value = 1726.3182
$ws.Range("J126").Value = 1749
$ws.Range("K126").Value = 5178.9546
$ws.Range("L126").Value = 5247
$ws.Range("M126").Value = -2708.9546
$ws.Range("N126").Value = -10187
